# Scheduled runner update: refresh market price / profit figures across all crafting-sheet tables.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Range("H19").Value2 = 2027.8889
$ws.Range("I19").Value2 = 2581.6924
$ws.Range("J19").Value2 = 1513.6428
$ws.Range("K19").Value2 = 2581.6924
$ws.Range("L19").Value2 = 1513.6428
$ws.Range("M19").Value2 = -2406.6924
$ws.Range("N19").Value2 = -1863.6428
# Row 41
$ws.Range("H41").Value2 = 1057.8572
$ws.Range("I41").Value2 = 1103.2222
$ws.Range("J41").Value2 = 976.2
$ws.Range("K41").Value2 = 1103.2222
$ws.Range("L41").Value2 = 976.2
$ws.Range("M41").Value2 = -663.2221999999999
$ws.Range("N41").Value2 = -1856.2
# Row 54
$ws.Range("H54").Value2 = 15999.857
$ws.Range("I54").Value2 = 15999.857
$ws.Range("K54").Value2 = 15999.857
$ws.Range("M54").Value2 = -15513.857
# Row 86
$ws.Range("H86").Value2 = 10778961
$ws.Range("I86").Value2 = 4168812.5
$ws.Range("K86").Value2 = 4168812.5
$ws.Range("M86").Value2 = -4167689.5
# Row 89
$ws.Range("H89").Value2 = 10778961
$ws.Range("I89").Value2 = 4168812.5
$ws.Range("K89").Value2 = 20844062.5
$ws.Range("M89").Value2 = -20838446.5
# Row 106
$ws.Range("H106").Value2 = 4337.615
$ws.Range("I106").Value2 = 4217.273
$ws.Range("J106").Value2 = 4999.5
$ws.Range("K106").Value2 = 4217.273
$ws.Range("L106").Value2 = 4999.5
$ws.Range("M106").Value2 = -3586.273
$ws.Range("N106").Value2 = -6261.5
# Row 113
$ws.Range("H113").Value2 = 2624
$ws.Range("I113").Value2 = 2648.8
$ws.Range("J113").Value2 = 2500
$ws.Range("K113").Value2 = 2648.8
$ws.Range("L113").Value2 = 2500
$ws.Range("M113").Value2 = 605.1999999999998
$ws.Range("N113").Value2 = -9008
# Row 116
$ws.Range("H116").Value2 = 3863.7334
$ws.Range("I116").Value2 = 3707.1428
$ws.Range("J116").Value2 = 4000.75
$ws.Range("K116").Value2 = 3707.1428
$ws.Range("L116").Value2 = 4000.75
$ws.Range("M116").Value2 = -265.1428000000001
$ws.Range("N116").Value2 = -10884.75
# Row 138
$ws.Range("H138").Value2 = 3928.2334
$ws.Range("I138").Value2 = 2160.7144
$ws.Range("K138").Value2 = 6482.1432
$ws.Range("M138").Value2 = -1342.1432

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value2 = 4790
$ws.Range("I2").Value2 = 3488.9285
$ws.Range("K2").Value2 = 3488.9285
$ws.Range("M2").Value2 = -3375.9285
# Row 32
$ws.Range("H32").Value2 = 940.9873700000001
$ws.Range("I32").Value2 = 742.04285
$ws.Range("K32").Value2 = 742.04285
$ws.Range("M32").Value2 = -455.04285
# Row 61
$ws.Range("H61").Value2 = 8923.706
$ws.Range("J61").Value2 = 6598.6
$ws.Range("L61").Value2 = 6598.6
$ws.Range("N61").Value2 = -7022.6
# Row 116
$ws.Range("H116").Value2 = 4790
$ws.Range("I116").Value2 = 3488.9285
$ws.Range("K116").Value2 = 3488.9285
$ws.Range("M116").Value2 = -1194.9285
# Row 122
$ws.Range("H122").Value2 = 2482.4146
$ws.Range("I122").Value2 = 1663.6428
$ws.Range("K122").Value2 = 4990.928400000001
$ws.Range("M122").Value2 = -2540.928400000001
# Row 132
$ws.Range("H132").Value2 = 2632.532
$ws.Range("I132").Value2 = 1934.0476
$ws.Range("J132").Value2 = 8499.799999999999
$ws.Range("K132").Value2 = 5802.142800000001
$ws.Range("L132").Value2 = 25499.4
$ws.Range("M132").Value2 = -3272.142800000001
$ws.Range("N132").Value2 = -30559.4
# Row 136
$ws.Range("H136").Value2 = 8923.706
$ws.Range("J136").Value2 = 6598.6
$ws.Range("L136").Value2 = 19795.8
$ws.Range("N136").Value2 = -24895.8

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value2 = 4790
$ws.Range("I3").Value2 = 3488.9285
$ws.Range("K3").Value2 = 3488.9285
$ws.Range("M3").Value2 = -3374.9285
# Row 107
$ws.Range("H107").Value2 = 3906.9
$ws.Range("I107").Value2 = 3906.9
$ws.Range("K107").Value2 = 3906.9
$ws.Range("M107").Value2 = -1986.9

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value2 = 5317.8237
$ws.Range("I31").Value2 = 12549.125
$ws.Range("J31").Value2 = 3092.8076
$ws.Range("K31").Value2 = 12549.125
$ws.Range("L31").Value2 = 3092.8076
$ws.Range("M31").Value2 = -12254.125
$ws.Range("N31").Value2 = -3682.8076
# Row 34
$ws.Range("H34").Value2 = 5317.8237
$ws.Range("I34").Value2 = 12549.125
$ws.Range("J34").Value2 = 3092.8076
$ws.Range("K34").Value2 = 12549.125
$ws.Range("L34").Value2 = 3092.8076
$ws.Range("M34").Value2 = -12347.125
$ws.Range("N34").Value2 = -3496.8076
# Row 134
$ws.Range("H134").Value2 = 1675.2307
$ws.Range("I134").Value2 = 1681.5834
$ws.Range("K134").Value2 = 5044.7502
$ws.Range("M134").Value2 = -2509.7502

$ws = $wb.Worksheets.Item("CUL")
# Row 129
$ws.Range("H129").Value2 = 1635.7693
$ws.Range("I129").Value2 = 428.14285
$ws.Range("J129").Value2 = 3044.6667
$ws.Range("K129").Value2 = 1284.42855
$ws.Range("L129").Value2 = 9134.000100000001
$ws.Range("M129").Value2 = 3715.57145
$ws.Range("N129").Value2 = -19134.0001
# Row 137
$ws.Range("H137").Value2 = 4222.231
$ws.Range("I137").Value2 = 1900
$ws.Range("J137").Value2 = 4644.4546
$ws.Range("K137").Value2 = 5700
$ws.Range("L137").Value2 = 13933.3638
$ws.Range("M137").Value2 = -600
$ws.Range("N137").Value2 = -24133.3638

$ws = $wb.Worksheets.Item("GSM")
# Row 33
$ws.Range("H33").Value2 = 47899.668
$ws.Range("I33").Value2 = 47800
$ws.Range("K33").Value2 = 47800
$ws.Range("M33").Value2 = -47548
# Row 41
$ws.Range("H41").Value2 = 4849.5
$ws.Range("I41").Value2 = 9998
$ws.Range("J41").Value2 = 3133.3333
$ws.Range("K41").Value2 = 9998
$ws.Range("L41").Value2 = 3133.3333
$ws.Range("M41").Value2 = -9643
$ws.Range("N41").Value2 = -3843.3333
# Row 132
$ws.Range("H132").Value2 = 12047.451
$ws.Range("I132").Value2 = 6346.48
$ws.Range("K132").Value2 = 19039.44
$ws.Range("M132").Value2 = -16509.44

$ws = $wb.Worksheets.Item("LTW")
# Row 61
$ws.Range("H61").Value2 = 2884.7144
$ws.Range("I61").Value2 = 2868.0908
$ws.Range("J61").Value2 = 2945.6667
$ws.Range("K61").Value2 = 2868.0908
$ws.Range("L61").Value2 = 2945.6667
$ws.Range("M61").Value2 = -2666.0908
$ws.Range("N61").Value2 = -3349.6667
# Row 68
$ws.Range("H68").Value2 = 3087.2
$ws.Range("I68").Value2 = 2919.4
$ws.Range("J68").Value2 = 3255
$ws.Range("K68").Value2 = 2919.4
$ws.Range("L68").Value2 = 3255
$ws.Range("M68").Value2 = -2170.4
$ws.Range("N68").Value2 = -4753
# Row 71
$ws.Range("H71").Value2 = 3087.2
$ws.Range("I71").Value2 = 2919.4
$ws.Range("J71").Value2 = 3255
$ws.Range("K71").Value2 = 14597
$ws.Range("L71").Value2 = 16275
$ws.Range("M71").Value2 = -10853
$ws.Range("N71").Value2 = -23763
# Row 82
$ws.Range("H82").Value2 = 2046.3334
$ws.Range("I82").Value2 = 1794.2858
$ws.Range("J82").Value2 = 2399.2
$ws.Range("K82").Value2 = 1794.2858
$ws.Range("L82").Value2 = 2399.2
$ws.Range("M82").Value2 = -1433.2858
$ws.Range("N82").Value2 = -3121.2
# Row 85
$ws.Range("H85").Value2 = 2046.3334
$ws.Range("I85").Value2 = 1794.2858
$ws.Range("J85").Value2 = 2399.2
$ws.Range("K85").Value2 = 1794.2858
$ws.Range("L85").Value2 = 2399.2
$ws.Range("M85").Value2 = -546.2858000000001
$ws.Range("N85").Value2 = -4895.2
# Row 113
$ws.Range("H113").Value2 = 2884.7144
$ws.Range("I113").Value2 = 2868.0908
$ws.Range("J113").Value2 = 2945.6667
$ws.Range("K113").Value2 = 2868.0908
$ws.Range("L113").Value2 = 2945.6667
$ws.Range("M113").Value2 = -698.0907999999999
$ws.Range("N113").Value2 = -7285.6667
# Row 132
$ws.Range("H132").Value2 = 3228.3044
$ws.Range("I132").Value2 = 3228.3044
$ws.Range("J132").Value2 = 0
$ws.Range("K132").Value2 = 9684.913199999999
$ws.Range("L132").Value2 = 0
$ws.Range("M132").Value2 = -7154.913199999999
$ws.Range("N132").ClearContents()
# Row 136
$ws.Range("H136").Value2 = 4896.091
$ws.Range("J136").Value2 = 14444
$ws.Range("L136").Value2 = 43332
$ws.Range("N136").Value2 = -48432

$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value2 = 3090.7778
$ws.Range("I132").Value2 = 1831.9259
$ws.Range("K132").Value2 = 5495.7777
$ws.Range("M132").Value2 = -2965.7777
# Row 136
$ws.Range("H136").Value2 = 3376
$ws.Range("I136").Value2 = 2720.2354
$ws.Range("K136").Value2 = 8160.706200000001
$ws.Range("M136").Value2 = -5610.706200000001
